$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 1233 (pushing the existing rows 1233:1308 down to 1236:1311)
$ws.Rows("1233:1235").Insert()

# Copy the date style (s="2", custom date/time number format) used by column D onto the new D cells
$ws.Range("D1236").Copy()
$ws.Range("D1233:D1235").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 1233 (Cebollín, Extra, fecha 44585)
$ws.Cells.Item(1233, 1).Value = 6
$ws.Cells.Item(1233, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1233, 3).Value = "Metropolitana"
$ws.Cells.Item(1233, 4).Value = 44585
$ws.Cells.Item(1233, 5).Value = 13
$ws.Cells.Item(1233, 6).Value = 100112037
$ws.Cells.Item(1233, 7).Value = "Cebollín"
$ws.Cells.Item(1233, 8).Value = "Sin especificar"
$ws.Cells.Item(1233, 9).Value = "Extra"
$ws.Cells.Item(1233, 10).Value = 810
$ws.Cells.Item(1233, 11).Value = 2400
$ws.Cells.Item(1233, 12).Value = 2500
$ws.Cells.Item(1233, 13).Value = 2454
$ws.Cells.Item(1233, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(1233, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1233, 16).Value = 68
$ws.Cells.Item(1233, 17).Value = 36
$ws.Cells.Item(1233, 18).Value = "Hortaliza"

# Fill in the new row 1234 (Cebollín, Primera, fecha 44585)
$ws.Cells.Item(1234, 1).Value = 6
$ws.Cells.Item(1234, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1234, 3).Value = "Metropolitana"
$ws.Cells.Item(1234, 4).Value = 44585
$ws.Cells.Item(1234, 5).Value = 13
$ws.Cells.Item(1234, 6).Value = 100112037
$ws.Cells.Item(1234, 7).Value = "Cebollín"
$ws.Cells.Item(1234, 8).Value = "Sin especificar"
$ws.Cells.Item(1234, 9).Value = "Primera"
$ws.Cells.Item(1234, 10).Value = 1010
$ws.Cells.Item(1234, 11).Value = 2000
$ws.Cells.Item(1234, 12).Value = 2200
$ws.Cells.Item(1234, 13).Value = 2103
$ws.Cells.Item(1234, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(1234, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1234, 16).Value = 58
$ws.Cells.Item(1234, 17).Value = 36
$ws.Cells.Item(1234, 18).Value = "Hortaliza"

# Fill in the new row 1235 (Cebollín, Segunda, fecha 44585)
$ws.Cells.Item(1235, 1).Value = 6
$ws.Cells.Item(1235, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1235, 3).Value = "Metropolitana"
$ws.Cells.Item(1235, 4).Value = 44585
$ws.Cells.Item(1235, 5).Value = 13
$ws.Cells.Item(1235, 6).Value = 100112037
$ws.Cells.Item(1235, 7).Value = "Cebollín"
$ws.Cells.Item(1235, 8).Value = "Sin especificar"
$ws.Cells.Item(1235, 9).Value = "Segunda"
$ws.Cells.Item(1235, 10).Value = 300
$ws.Cells.Item(1235, 11).Value = 1800
$ws.Cells.Item(1235, 12).Value = 1800
$ws.Cells.Item(1235, 13).Value = 1800
$ws.Cells.Item(1235, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(1235, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1235, 16).Value = 50
$ws.Cells.Item(1235, 17).Value = 36
$ws.Cells.Item(1235, 18).Value = "Hortaliza"

Write-Host "Done"
